# Backlog / Requisitos cleanup
# - Renumber/retext the "RF#" requirement rows into a numeric ID column (1-8)
#   with corrected requirement text (accents fixed, wording tweaks).
# - Center the new numeric ID column (and its header).
# - Move the active selection to C12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (ID / Requisito / Classificação / Requisito) ---------------
# Only the "ID" header cell (B10) gains centered alignment in the new layout.
$ws.Range("B10").HorizontalAlignment = -4108   # xlCenter

# --- Data rows 11-18: numeric ID in col B (centered), updated text ---------

# Row 11
$ws.Range("B11").Value = 1
$ws.Range("B11").HorizontalAlignment = -4108
$ws.Range("C11").Value = "O software deve permitir que o administrador efetue a manutenção (consulta, inclusão, exclusão e alteração) dos dados."
$ws.Range("D11").Value = "Essencial"
$ws.Range("E11").Value = "Funcional"

# Row 12
$ws.Range("B12").Value = 2
$ws.Range("B12").HorizontalAlignment = -4108
$ws.Range("C12").Value = "O software deve permitir que o usuário solicite a consulta das vagas."
$ws.Range("D12").Value = "Essencial"
$ws.Range("E12").Value = "Funcional"

# Row 13
$ws.Range("B13").Value = 3
$ws.Range("B13").HorizontalAlignment = -4108
$ws.Range("C13").Value = "O software deve permitir que o usuário gere relatorios atraves dos dados obtidos."
$ws.Range("D13").Value = "Essencial"
$ws.Range("E13").Value = "Funcional"

# Row 14
$ws.Range("B14").Value = 4
$ws.Range("B14").HorizontalAlignment = -4108
$ws.Range("C14").Value = "O software deve permitir que o administrador gere relatórios através dos dados obtidos."
$ws.Range("D14").Value = "Importante"
$ws.Range("E14").Value = "Funcional"

# Row 15
$ws.Range("B15").Value = 5
$ws.Range("B15").HorizontalAlignment = -4108
$ws.Range("C15").Value = "O software deve gerar gráficos através dos dados obtidos."
$ws.Range("D15").Value = "Essencial"
$ws.Range("E15").Value = "Funcional"

# Row 16
$ws.Range("B16").Value = 6
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("C16").Value = "O software deve receber dados pelo sensor do arduíno."
$ws.Range("D16").Value = "Essencial"
$ws.Range("E16").Value = "Funcional"

# Row 17
$ws.Range("B17").Value = 7
$ws.Range("B17").HorizontalAlignment = -4108
$ws.Range("C17").Value = "O arduíno deve ter sensor de obstáculo."
$ws.Range("D17").Value = "Essencial"
$ws.Range("E17").Value = "Funcional"

# Row 18
$ws.Range("B18").Value = 8
$ws.Range("B18").HorizontalAlignment = -4108
$ws.Range("C18").Value = "O software deve permitir que o usuário comum solicite a consulta de vagas disponíveis no momento."
$ws.Range("D18").Value = "Desejável"
$ws.Range("E18").Value = "Funcional"

# --- Selection moves to C12 -------------------------------------------------
$ws.Range("C12").Select()
